$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Recorded By" values in column G (data rows 2 through 157),
# leaving the "Recorded By" header in G1 untouched.
$ws.Range("G2:G157").ClearContents()

# Narrow column G from its original width (31 characters) down to 13
# characters. 12.1666666667 is the ColumnWidth input that this runtime's
# character->pixel->character rounding maps back to a stored width of
# exactly 13.
$ws.Columns.Item(7).ColumnWidth = 12.1666666667
